$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1069.4166
$ws.Range("I88").Value = 640.5
$ws.Range("J88").Value = 1498.3334
$ws.Range("K88").Value = 640.5
$ws.Range("L88").Value = 1498.3334
$ws.Range("M88").Value = -234.5
$ws.Range("N88").Value = -2310.3334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 1069.4166
$ws.Range("I91").Value = 640.5
$ws.Range("J91").Value = 1498.3334
$ws.Range("K91").Value = 640.5
$ws.Range("L91").Value = 1498.3334
$ws.Range("M91").Value = 763.5
$ws.Range("N91").Value = -4306.3334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 1500
$ws.Range("I113").Value = 1500
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 1754
$ws.Range("N113").Value = -8008

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 1821.2222
$ws.Range("I116").Value = 1650
$ws.Range("J116").Value = 2163.6667
$ws.Range("K116").Value = 1650
$ws.Range("L116").Value = 2163.6667
$ws.Range("M116").Value = 1792

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 866622.9399999999
$ws.Range("I2").Value = 1416.591
$ws.Range("J2").Value = 2452834.5
$ws.Range("K2").Value = 1416.591
$ws.Range("L2").Value = 2452834.5
$ws.Range("M2").Value = -1303.591

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 15794296
$ws.Range("I74").Value = 21429574
$ws.Range("J74").Value = 15514.5
$ws.Range("K74").Value = 21429574
$ws.Range("L74").Value = 15514.5
$ws.Range("M74").Value = -21428700
$ws.Range("N74").Value = -17262.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 15794296
$ws.Range("I77").Value = 21429574
$ws.Range("J77").Value = 15514.5
$ws.Range("K77").Value = 107147870
$ws.Range("L77").Value = 77572.5
$ws.Range("M77").Value = -107143502
$ws.Range("N77").Value = -86308.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 548.6
$ws.Range("I97").Value = 469.08334
$ws.Range("J97").Value = 866.6667
$ws.Range("K97").Value = 469.08334
$ws.Range("L97").Value = 866.6667
$ws.Range("M97").Value = 26.91665999999998
$ws.Range("N97").Value = -1858.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 866622.9399999999
$ws.Range("I116").Value = 1416.591
$ws.Range("J116").Value = 2452834.5
$ws.Range("K116").Value = 1416.591
$ws.Range("L116").Value = 2452834.5
$ws.Range("M116").Value = 877.4090000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 889479.75
$ws.Range("I132").Value = 1570175.8
$ws.Range("J132").Value = 57517.945
$ws.Range("K132").Value = 4710527.4
$ws.Range("L132").Value = 172553.835
$ws.Range("M132").Value = -4707997.4
$ws.Range("N132").Value = -177613.835

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H140").Value = 36660.75
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 36660.75
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 36660.75
$ws.Range("N140").Value = -47020.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 866622.9399999999
$ws.Range("I3").Value = 1416.591
$ws.Range("J3").Value = 2452834.5
$ws.Range("K3").Value = 1416.591
$ws.Range("L3").Value = 2452834.5
$ws.Range("M3").Value = -1302.591

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2962.3333
$ws.Range("I86").Value = 1825.3334
$ws.Range("J86").Value = 4667.8335
$ws.Range("K86").Value = 1825.3334
$ws.Range("L86").Value = 4667.8335
$ws.Range("M86").Value = -702.3334
$ws.Range("N86").Value = -6913.8335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2962.3333
$ws.Range("I89").Value = 1825.3334
$ws.Range("J89").Value = 4667.8335
$ws.Range("K89").Value = 9126.666999999999
$ws.Range("L89").Value = 23339.1675
$ws.Range("M89").Value = -3510.666999999999
$ws.Range("N89").Value = -34571.1675

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 715.2778
$ws.Range("I94").Value = 715.2778
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 715.2778
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -264.2778
$ws.Range("N94").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 62613.5
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 62613.5
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 62613.5
$ws.Range("N140").Value = -72973.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2744.7856
$ws.Range("I62").Value = 2610.8333
$ws.Range("J62").Value = 2845.25
$ws.Range("K62").Value = 2610.8333
$ws.Range("L62").Value = 2845.25
$ws.Range("M62").Value = -1986.8333
$ws.Range("N62").Value = -4093.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 2744.7856
$ws.Range("I65").Value = 2610.8333
$ws.Range("J65").Value = 2845.25
$ws.Range("K65").Value = 13054.1665
$ws.Range("L65").Value = 14226.25
$ws.Range("M65").Value = -9934.166499999999
$ws.Range("N65").Value = -20466.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1985648.8
$ws.Range("I107").Value = 5211272
$ws.Range("J107").Value = 649.6923
$ws.Range("K107").Value = 5211272
$ws.Range("L107").Value = 649.6923
$ws.Range("M107").Value = -5209352
$ws.Range("N107").Value = -4489.6923

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2660
$ws.Range("I134").Value = 2528.25
$ws.Range("J134").Value = 3081.6
$ws.Range("K134").Value = 7584.75
$ws.Range("L134").Value = 9244.799999999999
$ws.Range("M134").Value = -5049.75
$ws.Range("N134").Value = -14314.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 52.5
$ws.Range("I8").Value = 52.5
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 157.5
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -18.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1894.3077
$ws.Range("I131").Value = 5029.8184
$ws.Range("J131").Value = 1463.175
$ws.Range("K131").Value = 15089.4552
$ws.Range("L131").Value = 4389.525
$ws.Range("M131").Value = -10049.4552
$ws.Range("N131").Value = -14469.525

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4601.3335
$ws.Range("I70").Value = 3904
$ws.Range("J70").Value = 4950
$ws.Range("K70").Value = 3904
$ws.Range("L70").Value = 4950
$ws.Range("M70").Value = -3634
$ws.Range("N70").Value = -5490

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4601.3335
$ws.Range("I73").Value = 3904
$ws.Range("J73").Value = 4950
$ws.Range("K73").Value = 3904
$ws.Range("L73").Value = 4950
$ws.Range("M73").Value = -2968
$ws.Range("N73").Value = -6822

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 61230.65
$ws.Range("I80").Value = 2645.5557
$ws.Range("J80").Value = 127138.875
$ws.Range("K80").Value = 2645.5557
$ws.Range("L80").Value = 127138.875
$ws.Range("M80").Value = -1647.5557
$ws.Range("N80").Value = -129134.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 61230.65
$ws.Range("I83").Value = 2645.5557
$ws.Range("J83").Value = 127138.875
$ws.Range("K83").Value = 13227.7785
$ws.Range("L83").Value = 635694.375
$ws.Range("M83").Value = -8235.7785
$ws.Range("N83").Value = -645678.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2161.0527
$ws.Range("I97").Value = 2124
$ws.Range("J97").Value = 2300
$ws.Range("K97").Value = 2124
$ws.Range("L97").Value = 2300
$ws.Range("M97").Value = -1628
$ws.Range("N97").Value = -3292

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2405.4285
$ws.Range("I126").Value = 1982.75
$ws.Range("J126").Value = 2969
$ws.Range("K126").Value = 5948.25
$ws.Range("L126").Value = 8907
$ws.Range("M126").Value = -3478.25
$ws.Range("N126").Value = -13847

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 16212.875
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 16212.875
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 48638.625
$ws.Range("N136").Value = -53738.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1409.1072
$ws.Range("I82").Value = 1231.4166
$ws.Range("J82").Value = 1542.375
$ws.Range("K82").Value = 1231.4166
$ws.Range("L82").Value = 1542.375
$ws.Range("M82").Value = -870.4166
$ws.Range("N82").Value = -2264.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1409.1072
$ws.Range("I85").Value = 1231.4166
$ws.Range("J85").Value = 1542.375
$ws.Range("K85").Value = 1231.4166
$ws.Range("L85").Value = 1542.375
$ws.Range("M85").Value = 16.58339999999998
$ws.Range("N85").Value = -4038.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1824.1428
$ws.Range("I93").Value = 1853.8
$ws.Range("J93").Value = 1750
$ws.Range("K93").Value = 1853.8
$ws.Range("L93").Value = 1750
$ws.Range("M93").Value = -605.8
$ws.Range("N93").Value = -4246

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2536.0557
$ws.Range("I122").Value = 2359.1538
$ws.Range("J122").Value = 2996
$ws.Range("K122").Value = 7077.4614
$ws.Range("L122").Value = 8988
$ws.Range("M122").Value = -4627.4614

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3896.1
$ws.Range("I81").Value = 2199.8
$ws.Range("J81").Value = 5592.4
$ws.Range("K81").Value = 4399.6
$ws.Range("L81").Value = 11184.8
$ws.Range("M81").Value = -3338.6
$ws.Range("N81").Value = -13306.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 3896.1
$ws.Range("I84").Value = 2199.8
$ws.Range("J84").Value = 5592.4
$ws.Range("K84").Value = 21998
$ws.Range("L84").Value = 55924
$ws.Range("M84").Value = -16694
$ws.Range("N84").Value = -66532

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 277.94736
$ws.Range("I107").Value = 231
$ws.Range("J107").Value = 342.5
$ws.Range("K107").Value = 693
$ws.Range("L107").Value = 1027.5
$ws.Range("M107").Value = 1227
$ws.Range("N107").Value = -4867.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 490.42856
$ws.Range("I113").Value = 397.29413
$ws.Range("J113").Value = 634.36365
$ws.Range("K113").Value = 1191.88239
$ws.Range("L113").Value = 1903.09095
$ws.Range("M113").Value = 978.11761
$ws.Range("N113").Value = -6243.09095
